# Daily attendance processing - 2025-10-26 14:43:10
# Reorders the list of emails in the "Recorded By" (column G) cell for the rows below.
# The set of recorded-by emails for each row is unchanged; only the order within the
# comma-separated list is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range("G4").Value = 'hend_mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Range("G9").Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G10").Value = 'wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G15").Value = 'nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G16").Value = 'Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range("G23").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G24").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G29").Value = 'ola.m.abdelfattah@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G34").Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range("G35").Value = 'hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G40").Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range("G41").Value = 'aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G54").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G55").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G60").Value = 'ola.m.abdelfattah@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G64").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G65").Value = 'rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range("G66").Value = 'hend_mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Range("G71").Value = 'mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G72").Value = 'wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G76").Value = 'nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G85").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G86").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G89").Value = 'Aya.alaa-eldein@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G91").Value = 'ola.m.abdelfattah@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G95").Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range("G96").Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G97").Value = 'asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G102").Value = 'mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G103").Value = 'aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G108").Value = 'nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G109").Value = 'Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range("G120").Value = 'nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G122").Value = 'Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G126").Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range("G127").Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G128").Value = 'asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G133").Value = 'mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G138").Value = 'nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G139").Value = 'Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G141").Value = 'nourhanmamdouh@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range("G151").Value = 'user@user.com, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Range("G153").Value = 'Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G154").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range("G157").Value = 'hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G158").Value = 'rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range("G159").Value = 'hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G164").Value = 'mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G170").Value = 'Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range("G172").Value = 'nourhanmamdouh@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range("G178").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G179").Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Range("G182").Value = 'nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G184").Value = 'maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G185").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
